$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -16.53598913478064
$ws.Range("C2").Value = -0.3090313680422698
$ws.Range("D2").Value = -16.53598913478064
$ws.Range("E2").Value = -16.53598913478064
$ws.Range("F2").Value = -16.53598913478064
$ws.Range("G2").Value = -16.53598913478064
$ws.Range("H2").Value = -16.53598913478064
$ws.Range("I2").Value = -16.53598913478064
$ws.Range("J2").Value = -16.53598913478064
$ws.Range("K2").Value = -16.53598913478064

$ws.Range("B3").Value = -16.53598913478064
$ws.Range("C3").Value = -16.53598913478064
$ws.Range("D3").Value = -16.53598913478064
$ws.Range("E3").Value = -16.53598913478064
$ws.Range("F3").Value = -16.53598913478064
$ws.Range("G3").Value = -16.53598913478064
$ws.Range("H3").Value = -16.53598913478064
$ws.Range("I3").Value = 0.6075236257491098
$ws.Range("J3").Value = -16.53598913478064
$ws.Range("K3").Value = -16.53598913478064

$ws.Range("B4").Value = -16.53598913478064
$ws.Range("C4").Value = -0.2722614614953611
$ws.Range("D4").Value = 0.3772392286187093
$ws.Range("E4").Value = -16.53598913478064
$ws.Range("F4").Value = 3.8557883252331
$ws.Range("G4").Value = -16.53598913478064
$ws.Range("H4").Value = -16.53598913478064
$ws.Range("I4").Value = -16.53598913478064
$ws.Range("J4").Value = 1.950771268112282
$ws.Range("K4").Value = -16.53598913478064

$ws.Range("B5").Value = -16.53598913478064
$ws.Range("C5").Value = 0.405711884235391
$ws.Range("D5").Value = -16.53598913478064
$ws.Range("E5").Value = -16.53598913478064
$ws.Range("F5").Value = -16.53598913478064
$ws.Range("G5").Value = 3.092750077940472
$ws.Range("H5").Value = -16.53598913478064
$ws.Range("I5").Value = -16.53598913478064
$ws.Range("J5").Value = -16.53598913478064
$ws.Range("K5").Value = -16.53598913478064

$ws.Range("B6").Value = -16.53598913478064
$ws.Range("C6").Value = -16.53598913478064
$ws.Range("D6").Value = -16.53598913478064
$ws.Range("E6").Value = -16.53598913478064
$ws.Range("F6").Value = -16.53598913478064
$ws.Range("G6").Value = -16.53598913478064
$ws.Range("H6").Value = -16.53598913478064
$ws.Range("I6").Value = -16.53598913478064
$ws.Range("J6").Value = -16.53598913478064
$ws.Range("K6").Value = -16.53598913478064

$ws.Range("B7").Value = 3.1784593146922
$ws.Range("C7").Value = -16.53598913478064
$ws.Range("D7").Value = -16.53598913478064
$ws.Range("E7").Value = -16.53598913478064
$ws.Range("F7").Value = -16.53598913478064
$ws.Range("G7").Value = -16.53598913478064
$ws.Range("H7").Value = -16.53598913478064
$ws.Range("I7").Value = -16.53598913478064
$ws.Range("J7").Value = -16.53598913478064
$ws.Range("K7").Value = -16.53598913478064

$ws.Range("B8").Value = -16.53598913478064
$ws.Range("C8").Value = -16.53598913478064
$ws.Range("D8").Value = -16.53598913478064
$ws.Range("E8").Value = 1.490820648408289
$ws.Range("F8").Value = -16.53598913478064
$ws.Range("G8").Value = -16.53598913478064
$ws.Range("H8").Value = -16.53598913478064
$ws.Range("I8").Value = -16.53598913478064
$ws.Range("J8").Value = -16.53598913478064
$ws.Range("K8").Value = -16.53598913478064

$ws.Range("B9").Value = 3.452386285914278
$ws.Range("C9").Value = -16.53598913478064
$ws.Range("D9").Value = -16.53598913478064
$ws.Range("E9").Value = -16.53598913478064
$ws.Range("F9").Value = -16.53598913478064
$ws.Range("G9").Value = -16.53598913478064
$ws.Range("H9").Value = -16.53598913478064
$ws.Range("I9").Value = -16.53598913478064
$ws.Range("J9").Value = -16.53598913478064
$ws.Range("K9").Value = -16.53598913478064

$ws.Range("B10").Value = -16.53598913478064
$ws.Range("C10").Value = -16.53598913478064
$ws.Range("D10").Value = -16.53598913478064
$ws.Range("E10").Value = -16.53598913478064
$ws.Range("F10").Value = -16.53598913478064
$ws.Range("G10").Value = -16.53598913478064
$ws.Range("H10").Value = -16.53598913478064
$ws.Range("I10").Value = 1.37219822325161
$ws.Range("J10").Value = -16.53598913478064
$ws.Range("K10").Value = 2.003900443374778

$ws.Range("B11").Value = -16.53598913478064
$ws.Range("C11").Value = -16.53598913478064
$ws.Range("D11").Value = -16.53598913478064
$ws.Range("E11").Value = 2.104202834985375
$ws.Range("F11").Value = -16.53598913478064
$ws.Range("G11").Value = 1.983620245405954
$ws.Range("H11").Value = -16.53598913478064
$ws.Range("I11").Value = -16.53598913478064
$ws.Range("J11").Value = -16.53598913478064
$ws.Range("K11").Value = 1.580766009051879

$ws.Range("B12").Value = -16.53598913478064
$ws.Range("C12").Value = -16.53598913478064
$ws.Range("D12").Value = -16.53598913478064
$ws.Range("E12").Value = -16.53598913478064
$ws.Range("F12").Value = -16.53598913478064
$ws.Range("G12").Value = -16.53598913478064
$ws.Range("H12").Value = -16.53598913478064
$ws.Range("I12").Value = -16.53598913478064
$ws.Range("J12").Value = -16.53598913478064
$ws.Range("K12").Value = -16.53598913478064

$ws.Range("B13").Value = -16.53598913478064
$ws.Range("C13").Value = -16.53598913478064
$ws.Range("D13").Value = -16.53598913478064
$ws.Range("E13").Value = 1.99904697790876
$ws.Range("F13").Value = -16.53598913478064
$ws.Range("G13").Value = -16.53598913478064
$ws.Range("H13").Value = -16.53598913478064
$ws.Range("I13").Value = -16.53598913478064
$ws.Range("J13").Value = 1.593701394004662
$ws.Range("K13").Value = 1.811774953057695

$ws.Range("B14").Value = -16.53598913478064
$ws.Range("C14").Value = -16.53598913478064
$ws.Range("D14").Value = 0.968030651028902
$ws.Range("E14").Value = -16.53598913478064
$ws.Range("F14").Value = -16.53598913478064
$ws.Range("G14").Value = -16.53598913478064
$ws.Range("H14").Value = -16.53598913478064
$ws.Range("I14").Value = -16.53598913478064
$ws.Range("J14").Value = -16.53598913478064
$ws.Range("K14").Value = 2.110052687691336

$ws.Range("B15").Value = -16.53598913478064
$ws.Range("C15").Value = -16.53598913478064
$ws.Range("D15").Value = -0.31430761915835
$ws.Range("E15").Value = -16.53598913478064
$ws.Range("F15").Value = -16.53598913478064
$ws.Range("G15").Value = -16.53598913478064
$ws.Range("H15").Value = -16.53598913478064
$ws.Range("I15").Value = -16.53598913478064
$ws.Range("J15").Value = -16.53598913478064
$ws.Range("K15").Value = -16.53598913478064

$ws.Range("B16").Value = -16.53598913478064
$ws.Range("C16").Value = -16.53598913478064
$ws.Range("D16").Value = -16.53598913478064
$ws.Range("E16").Value = -16.53598913478064
$ws.Range("F16").Value = -16.53598913478064
$ws.Range("G16").Value = -16.53598913478064
$ws.Range("H16").Value = -16.53598913478064
$ws.Range("I16").Value = -16.53598913478064
$ws.Range("J16").Value = 1.951478026966471
$ws.Range("K16").Value = -16.53598913478064

$ws.Range("B17").Value = -16.53598913478064
$ws.Range("C17").Value = 0.6412065220514327
$ws.Range("D17").Value = -0.1045405004558063
$ws.Range("E17").Value = -16.53598913478064
$ws.Range("F17").Value = -16.53598913478064
$ws.Range("G17").Value = -16.53598913478064
$ws.Range("H17").Value = -16.53598913478064
$ws.Range("I17").Value = 1.012358679456854
$ws.Range("J17").Value = 2.311401654043482
$ws.Range("K17").Value = -16.53598913478064

$ws.Range("B18").Value = -16.53598913478064
$ws.Range("C18").Value = -16.53598913478064
$ws.Range("D18").Value = -16.53598913478064
$ws.Range("E18").Value = -16.53598913478064
$ws.Range("F18").Value = -16.53598913478064
$ws.Range("G18").Value = -16.53598913478064
$ws.Range("H18").Value = -16.53598913478064
$ws.Range("I18").Value = 0.3626712930862772
$ws.Range("J18").Value = 2.099108979740021
$ws.Range("K18").Value = -16.53598913478064

$ws.Range("B19").Value = -16.53598913478064
$ws.Range("C19").Value = -16.53598913478064
$ws.Range("D19").Value = 2.993946545671792
$ws.Range("E19").Value = -16.53598913478064
$ws.Range("F19").Value = -16.53598913478064
$ws.Range("G19").Value = -16.53598913478064
$ws.Range("H19").Value = -16.53598913478064
$ws.Range("I19").Value = 2.155130729942673
$ws.Range("J19").Value = -16.53598913478064
$ws.Range("K19").Value = -16.53598913478064

$ws.Range("B20").Value = -16.53598913478064
$ws.Range("C20").Value = 3.227750047929727
$ws.Range("D20").Value = 2.816392727909827
$ws.Range("E20").Value = -16.53598913478064
$ws.Range("F20").Value = 2.465146879196011
$ws.Range("G20").Value = -16.53598913478064
$ws.Range("H20").Value = 4.321913671356848
$ws.Range("I20").Value = 3.023337537621007
$ws.Range("J20").Value = -16.53598913478064
$ws.Range("K20").Value = 2.370141051287527

$ws.Range("B21").Value = -16.53598913478064
$ws.Range("C21").Value = 2.611695276966768
$ws.Range("D21").Value = -16.53598913478064
$ws.Range("E21").Value = 3.152570241389514
$ws.Range("F21").Value = -16.53598913478064
$ws.Range("G21").Value = 2.909537022355233
$ws.Range("H21").Value = -16.53598913478064
$ws.Range("I21").Value = -16.53598913478064
$ws.Range("J21").Value = -16.53598913478064
$ws.Range("K21").Value = -16.53598913478064
